$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 286; this shifts the existing rows
# 286-356 down to 288-358 (and brings along cell styles, e.g. the date
# format on column D).
$ws.Rows("286:287").Insert()

# New row 286: Cebollín, Primera, Provincia de Quillota, date 44543
$ws.Cells.Item(286, 1).Value = 3
$ws.Cells.Item(286, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(286, 3).Value = "Coquimbo"
$ws.Cells.Item(286, 4).Value = 44543
$ws.Cells.Item(286, 5).Value = 5
$ws.Cells.Item(286, 6).Value = 100112037
$ws.Cells.Item(286, 7).Value = "Cebollín"
$ws.Cells.Item(286, 8).Value = "Sin especificar"
$ws.Cells.Item(286, 9).Value = "Primera"
$ws.Cells.Item(286, 10).Value = 160
$ws.Cells.Item(286, 11).Value = 3000
$ws.Cells.Item(286, 12).Value = 3000
$ws.Cells.Item(286, 13).Value = 3000
$ws.Cells.Item(286, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(286, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(286, 16).Value = 83
$ws.Cells.Item(286, 17).Value = 36
$ws.Cells.Item(286, 18).Value = "Hortaliza"

# New row 287: Cebollín, Segunda, Provincia de Quillota, date 44543
$ws.Cells.Item(287, 1).Value = 3
$ws.Cells.Item(287, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(287, 3).Value = "Coquimbo"
$ws.Cells.Item(287, 4).Value = 44543
$ws.Cells.Item(287, 5).Value = 5
$ws.Cells.Item(287, 6).Value = 100112037
$ws.Cells.Item(287, 7).Value = "Cebollín"
$ws.Cells.Item(287, 8).Value = "Sin especificar"
$ws.Cells.Item(287, 9).Value = "Segunda"
$ws.Cells.Item(287, 10).Value = 120
$ws.Cells.Item(287, 11).Value = 2000
$ws.Cells.Item(287, 12).Value = 2000
$ws.Cells.Item(287, 13).Value = 2000
$ws.Cells.Item(287, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(287, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(287, 16).Value = 56
$ws.Cells.Item(287, 17).Value = 36
$ws.Cells.Item(287, 18).Value = "Hortaliza"

# Row 343 (old row 341, shifted down by 2) keeps its original L value
# (3000) in the source commit even though every other field for that
# row follows the "shift by two" pattern — restore it explicitly in
# case the engine's row-insert/shift altered it.
$ws.Cells.Item(343, 12).Value = 3000
